$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "fgInwardQty" and "fgInwardUnit" columns (old G:H) and the
# "fgUnitQtyforFG" column (old L), shifting all subsequent columns left.
# Old layout (A..P): Date, productName, productCode, itemCode,
#   OutwardRawMaterial, OutwardUnit, fgInwardQty, fgInwardUnit, fgUnitQty,
#   fgUnitWt, fgTotalQty, fgUnitQtyforFG, finishedGoodsUnit, scrapQty,
#   scrapUnit, finishedBy
# New layout (A..M): Date, productName, productCode, itemCode,
#   OutwardRawMaterial, OutwardUnit, fgUnitQty, fgUnitWt, fgTotalQty,
#   finishedGoodsUnit, scrapQty, scrapUnit, finishedBy

$ws.Range("G1:H1").EntireColumn.Delete() | Out-Null
$ws.Range("J1").EntireColumn.Delete() | Out-Null

$ws.Range("A1").Select()

$wb.Save()
